$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "New_Registration"

# 2. Add the "Used" value in G2 (new shared string)
$ws.Range("G2").Value = "Used"

# 3. Header row (A1:G1) formatting: bold font + yellow fill.
#    Build the combined style on A1 first, then copy/paste the format onto
#    the rest of the header row so only ONE new cell style (bold+yellow)
#    gets added to the style table (matches target: cellXfs count 2 -> 3).
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.Interior.Color = 65535
$a1.Copy()
$headerRange = $ws.Range("A1:G1")
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Column G width / best-fit
$ws.Columns.Item(7).AutoFit()

# 5. Move the active selection to I2 (matches sheetView selection change)
$ws.Range("I2").Select() | Out-Null

# 6. Page setup - portrait orientation
$ws.PageSetup.Orientation = 1
